# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# Updates First_Noticeable_Increase_Index (C), First_Noticeable_Increase_Cumulative_Value (E)
# and Pulse_Width (G) on each of the Step3_DataPts_* sheets for rows 2-5.

$wb = $excel.ActiveWorkbook

# Column C (First_Noticeable_Increase_Index) and Column E
# (First_Noticeable_Increase_Cumulative_Value) are identical across all four
# Step3_DataPts_* sheets; Column G (Pulse_Width) differs per-sheet because it
# depends on the sheet-specific Point_Exceeds_Index (column D).
$commonUpdates = @{
    2 = @{ C = 58; E = 0.000331612557837233 }
    3 = @{ C = 59; E = 0.0007471044468722725 }
    4 = @{ C = 58; E = 0.001033359022648661 }
    5 = @{ C = 58; E = 0.0008824809586660525 }
}

$pulseWidthUpdates = @{
    "Step3_DataPts_0.5" = @{ 2 = 47; 3 = 49; 4 = 49; 5 = 49 }
    "Step3_DataPts_0.7" = @{ 2 = 69; 3 = 71; 4 = 71; 5 = 70 }
    "Step3_DataPts_0.8" = @{ 2 = 84; 3 = 86; 4 = 85; 5 = 85 }
    "Step3_DataPts_0.9" = @{ 2 = 109; 3 = 110; 4 = 110; 5 = 110 }
}

foreach ($sheetName in $pulseWidthUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 3).Value = $commonUpdates[$row].C
        $ws.Cells.Item($row, 5).Value = $commonUpdates[$row].E
        $ws.Cells.Item($row, 7).Value = $pulseWidthUpdates[$sheetName][$row]
    }
}
